# "updated main GSC export data"
#
# The Coverage.xlsx "Chart" sheet is a rolling daily export: row 2 always
# holds the oldest date. This refresh drops that oldest day (2025-10-22)
# and every later row shifts up by one, so the newest day already present
# (2026-01-17) simply becomes the new last row instead of a row being
# appended. Deleting the sheet row (rather than rewriting each cell) lets
# Excel itself shift the remaining data/shared-strings up, which is what
# produced this export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the oldest date's row; everything below shifts up one row.
$ws.Rows.Item(2).Delete()
